$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" "91.381.61"
Set-TextCell "E2" "  +3.55%  "

Set-TextCell "D3" "3.118.45"
Set-TextCell "E3" "  +1.67%  "

Set-TextCell "E4" "  +0.07%  "

Set-TextCell "D5" "219.79"
Set-TextCell "E5" "  +4.41%  "

Set-TextCell "D6" "623.32"
Set-TextCell "E6" "  +0.68%  "

Set-TextCell "E7" "  +2.72%  "

Set-TextCell "D8" "0.973"
Set-TextCell "E8" "  +21.38%  "

Set-TextCell "E9" "  +0.07%  "

Set-TextCell "D10" "3.116.84"
Set-TextCell "E10" "  +1.80%  "

Set-TextCell "D11" "0.717"
Set-TextCell "E11" "  +21.32%  "

Set-TextCell "E12" "  +5.42%  "

Set-TextCell "E13" "  +7.01%  "

Set-TextCell "D14" "34.56"
Set-TextCell "E14" "  +8.16%  "

Set-TextCell "D15" "91.187.63"
Set-TextCell "E15" "  +3.60%  "

Set-TextCell "E16" "  +1.98%  "

Set-TextCell "D17" "3.701.19"
Set-TextCell "E17" "  +1.99%  "

Set-TextCell "D18" "3.097.49"
Set-TextCell "E18" "  +1.25%  "

Set-TextCell "D19" "3.73"
Set-TextCell "E19" "  +13.94%  "

Set-TextCell "E20" "  +9.77%  "

Set-TextCell "D21" "14.12"
Set-TextCell "E21" "  +6.29%  "

Set-TextCell "D22" "437.91"
Set-TextCell "E22" "  +4.14%  "

Set-TextCell "D23" "8.79"
Set-TextCell "E23" "  +8.01%  "

Set-TextCell "E24" "  +6.23%  "

Set-TextCell "D25" "6.18"
Set-TextCell "E25" "  +13.08%  "

Set-TextCell "D28" "3.290.39"
Set-TextCell "E28" "  +1.93%  "

Set-TextCell "E29" "  -0.26%  "

Set-TextCell "E30" "  +0.89%  "

Set-TextCell "D31" "9.13"
Set-TextCell "E31" "  +13.73%  "

Set-TextCell "E32" "  -8.50%  "

Set-TextCell "D33" "526.54"
Set-TextCell "E33" "  +3.62%  "

Set-TextCell "D34" "3.75"
Set-TextCell "E34" "  +4.55%  "

Set-TextCell "D35" "7.10"
Set-TextCell "E35" "  +5.31%  "

Set-TextCell "D36" "0.142"
Set-TextCell "E36" "  +9.37%  "

Set-TextCell "D37" "23.79"
Set-TextCell "E37" "  +6.96%  "

Set-TextCell "E38" "  +3.61%  "

Set-TextCell "D39" "1.28"
Set-TextCell "E39" "  +2.99%  "

Set-TextCell "D40" "22.30"
Set-TextCell "E40" "  +0.48%  "

Set-TextCell "E43" "  +14.80%  "

Set-TextCell "E44" "  +0.00%  "

Set-TextCell "D45" "0.379"
Set-TextCell "E45" "  +5.22%  "

Set-TextCell "D46" "1.91"
Set-TextCell "E46" "  +5.65%  "

Set-TextCell "D47" "146.90"
Set-TextCell "E47" "  -0.54%  "

Set-TextCell "D48" "44.07"
Set-TextCell "E48" "  +1.71%  "

Set-TextCell "E49" "  +9.55%  "

Set-TextCell "D50" "167.40"
Set-TextCell "E50" "  +7.01%  "

Set-TextCell "D51" "0.000258"
Set-TextCell "E51" "  +20.20%  "

Set-TextCell "B26" "Litecoin"
Set-TextCell "C26" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D26" "87.38"
Set-TextCell "E26" "  +6.80%  "

Set-TextCell "B27" "Aptos"
Set-TextCell "C27" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D27" "12.22"
Set-TextCell "E27" "  +4.47%  "

Set-TextCell "B41" "Hedera"
Set-TextCell "C41" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D41" "0.0858"
Set-TextCell "E41" "  +26.66%  "

Set-TextCell "B42" "FirstDigitalUSD"
Set-TextCell "C42" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D42" "1.00"
Set-TextCell "E42" "  +0.08%  "
